# Apply cryptos list update (prices + 1h volume %) per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.737.88"
$ws.Range("E2").Value = "  +0.85%  "
$ws.Range("D3").Value = "3.424.65"
$ws.Range("E3").Value = "  +1.17%  "
$ws.Range("E4").Value = "  +0.08%  "
$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.20"
$ws.Range("D5").Style = $style
$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.77"
$ws.Range("D6").Style = $style
$ws.Range("D7").Value = "3.418.49"
$ws.Range("E7").Value = "  +1.18%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("E10").Value = "  +3.21%  "
$ws.Range("E11").Value = "  -1.14%  "
$style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "49.01"
$ws.Range("D12").Style = $style
$ws.Range("E12").Value = "  +0.82%  "
$ws.Range("E13").Value = "  +0.60%  "
$style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "690.45"
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = "  +2.00%  "
$ws.Range("D15").Value = "3.973.51"
$ws.Range("E15").Value = "  +1.12%  "
$style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.62"
$ws.Range("D16").Style = $style
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("D17").Value = "69.766.16"
$ws.Range("D18").Value = "3.424.40"
$ws.Range("E18").Value = "  +1.08%  "
$style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.121"
$ws.Range("D19").Style = $style
$ws.Range("E19").Value = "  +1.11%  "
$ws.Range("E20").Value = "  +0.04%  "
$style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.40"
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = "  +0.62%  "
$ws.Range("E22").Value = "  -0.50%  "
$style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.47"
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = "  +1.07%  "
$style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.93"
$ws.Range("D24").Style = $style
$ws.Range("E24").Value = "  -1.02%  "
$style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "100.80"
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = "  -2.52%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("E27").Value = "  -2.44%  "
$style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.59"
$ws.Range("D28").Style = $style
$ws.Range("E28").Value = "  -0.02%  "
$style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.47"
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = "  -1.98%  "
$style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.75"
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = "  +0.42%  "
$style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.17"
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = "  +2.53%  "
$style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "572.61"
$ws.Range("D32").Style = $style
$ws.Range("E32").Value = "  +3.33%  "
$style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.73"
$ws.Range("D33").Style = $style
$ws.Range("E33").Value = "  +1.21%  "
$ws.Range("E34").Value = "  -1.55%  "
$style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "58.34"
$ws.Range("D35").Style = $style
$ws.Range("E35").Value = "  +0.56%  "
$ws.Range("E36").Value = "  -2.56%  "
$ws.Range("E37").Value = "  +0.11%  "
$ws.Range("D38").Value = "3.599.49"
$ws.Range("E38").Value = "  -2.47%  "
$ws.Range("E39").Value = "  -0.24%  "
$style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "35.11"
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = "  +0.19%  "
$ws.Range("D41").Value = "0.0₃0740"
$ws.Range("E41").Value = "  +5.44%  "
$ws.Range("E42").Value = "  +0.48%  "
$ws.Range("E43").Value = "  +0.07%  "
$style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0420"
$ws.Range("D44").Style = $style
$ws.Range("E44").Value = "  -0.52%  "
$ws.Range("E45").Value = "  -1.56%  "
$ws.Range("E46").Value = "  +4.47%  "
$style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.66"
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("E48").Value = "  -0.40%  "
$ws.Range("E49").Value = "  -0.30%  "
$style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "133.36"
$ws.Range("D50").Style = $style
$ws.Range("E50").Value = "  +1.36%  "
$style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.65"
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = "  +2.15%  "
